# Fruta / hortaliza, semanal
# Update D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) columns per row,
# re-shuffling the weekly records according to the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44537; J = 88;  K = 2000; L = 2200; M = 2091; P = 697  }
    3  = @{ D = 44187; J = 65;  K = 3000; L = 3000; M = 3000; P = 1000 }
    4  = @{ D = 44223; J = 80;  K = 2500; L = 3000; M = 2781; P = 927  }
    5  = @{ D = 44225; J = 56;  K = 3000; L = 3000; M = 3000; P = 1000 }
    6  = @{ D = 44260; J = 60;  K = 3500; L = 3500; M = 3500; P = 1167 }
    7  = @{ D = 44557; J = 104; K = 2000; L = 2500; M = 2260; P = 753  }
    8  = @{ D = 44389; J = 81;  K = 2800; L = 3000; M = 2889; P = 963  }
    9  = @{ D = 44804; J = 85;  K = 3000; L = 3000; M = 3000; P = 1000 }
    10 = @{ D = 44390; J = 50;  K = 3000; L = 3000; M = 3000; P = 1000 }
    12 = @{ D = 44224; J = 67;  K = 3000; L = 3000; M = 3000; P = 1000 }
    13 = @{ D = 44165; J = 68;  K = 3000; L = 3000; M = 3000; P = 1000 }
    14 = @{ D = 44179; J = 78;  K = 3000; L = 3000; M = 3000; P = 1000 }
    15 = @{ D = 44291; J = 45;  K = 3000; L = 3000; M = 3000; P = 1000 }
    16 = @{ D = 44292; J = 40;  K = 3000; L = 3000; M = 3000; P = 1000 }
    17 = @{ D = 44536; J = 125; K = 2200; L = 2200; M = 2200; P = 733  }
    18 = @{ D = 44756; J = 104; K = 2800; L = 3000; M = 2904; P = 968  }
    19 = @{ D = 44222; J = 45;  K = 3000; L = 3000; M = 3000; P = 1000 }
    20 = @{ D = 44669; J = 92;  K = 2500; L = 3000; M = 2755; P = 918  }
    21 = @{ D = 44166; J = 45;  K = 2500; L = 2500; M = 2500; P = 833  }
    23 = @{ D = 44340; J = 54;  K = 3000; L = 3000; M = 3000; P = 1000 }
    24 = @{ D = 44242; J = 95;  K = 2500; L = 3000; M = 2737; P = 912  }
    25 = @{ D = 44221; J = 50;  K = 2500; L = 2500; M = 2500; P = 833  }
    26 = @{ D = 44559; J = 68;  K = 2000; L = 2000; M = 2000; P = 667  }
    27 = @{ D = 44845; J = 80;  K = 2500; L = 2500; M = 2500; P = 833  }
    28 = @{ D = 44627; J = 78;  K = 3500; L = 3500; M = 3500; P = 1167 }
    29 = @{ D = 44574; J = 50;  K = 3000; L = 3000; M = 3000; P = 1000 }
    30 = @{ D = 44193; J = 70;  K = 3000; L = 3000; M = 3000; P = 1000 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value2  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value2 = $vals.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value2 = $vals.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $vals.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value2 = $vals.P   # P - Precio $/Kg
}
